$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I61").Value = 67.23999999999999
$ws.Range("L69").Value = 68.73999999999999
$ws.Range("I71").Value = 72.37
$ws.Range("I96").Value = 80.23999999999999
$ws.Range("I98").Value = 84.12
$ws.Range("L100").Value = 80.48999999999999
$ws.Range("L103").Value = 81.73999999999999
$ws.Range("L104").Value = 82.12
$ws.Range("L105").Value = 82.48999999999999
$ws.Range("L112").Value = 85.12
$ws.Range("I114").Value = 88.37
$ws.Range("L120").Value = 88.98999999999999
$ws.Range("I126").Value = 99.73999999999999
$ws.Range("H132").Value = 935495
$ws.Range("I132").Value = 100.71
$ws.Range("M132").Value = 14067042
$ws.Range("M133").Value = 14469495
$ws.Range("M134").Value = 14605130
$ws.Range("M135").Value = 14690567
$ws.Range("M136").Value = 14867026
$ws.Range("M137").Value = 14861544
$ws.Range("M138").Value = 14670936
$ws.Range("M139").Value = 14467161
$ws.Range("D140").Value = 823711
$ws.Range("M140").Value = 14131941
$ws.Range("N140").Value = 8935547
$ws.Range("M141").Value = 13894237
$ws.Range("N141").Value = 8782062
$ws.Range("M142").Value = 13750144
$ws.Range("N142").Value = 8776330
$ws.Range("M143").Value = 13622013
$ws.Range("N143").Value = 8836772
$ws.Range("N144").Value = 9042580
$ws.Range("N145").Value = 9395200
$ws.Range("N146").Value = 9716064
$ws.Range("N147").Value = 9916895
$ws.Range("N148").Value = 10176528
$ws.Range("N149").Value = 10402109
$ws.Range("N150").Value = 10586509
$ws.Range("N151").Value = 10834422
$ws.Range("D186").Value = 590485
$ws.Range("G186").Value = 9818740
$ws.Range("H186").Value = 1208967
$ws.Range("M186").Value = 13889756
$ws.Range("N186").Value = 8211290
$ws.Range("M187").Value = 13965989
$ws.Range("N187").Value = 8231281
$ws.Range("M188").Value = 14177222
$ws.Range("N188").Value = 8311067
$ws.Range("M189").Value = 14260530
$ws.Range("N189").Value = 8315427
$ws.Range("L190").Value = 126.78
$ws.Range("M190").Value = 14422857
$ws.Range("N190").Value = 8371008
$ws.Range("D191").Value = 573551
$ws.Range("E191").Value = 207672
$ws.Range("G191").Value = 8670501
$ws.Range("H191").Value = 1073073
$ws.Range("M191").Value = 14528027
$ws.Range("N191").Value = 8368533
$ws.Range("O191").Value = 2551085
$ws.Range("M192").Value = 14532802
$ws.Range("N192").Value = 8290107
$ws.Range("O192").Value = 2547707
$ws.Range("M193").Value = 14761325
$ws.Range("N193").Value = 8365199
$ws.Range("O193").Value = 2572355
$ws.Range("F194").Value = 14729039
$ws.Range("M194").Value = 14844006
$ws.Range("N194").Value = 8388086
$ws.Range("O194").Value = 2570155
$ws.Range("M195").Value = 14888937
$ws.Range("N195").Value = 8325558
$ws.Range("O195").Value = 2571044
$ws.Range("F196").Value = 13097731
$ws.Range("M196").Value = 14967619
$ws.Range("N196").Value = 8320054
$ws.Range("O196").Value = 2563118
$ws.Range("M197").Value = 15089763
$ws.Range("N197").Value = 8341518
$ws.Range("O197").Value = 2566573
$ws.Range("M198").Value = 15141521
$ws.Range("N198").Value = 8330311
$ws.Range("O198").Value = 2567374
$ws.Range("F199").Value = 10187313
$ws.Range("M199").Value = 15309065
$ws.Range("N199").Value = 8356559
$ws.Range("O199").Value = 2587926
$ws.Range("M200").Value = 15440675
$ws.Range("N200").Value = 8350578
$ws.Range("O200").Value = 2599978
$ws.Range("M201").Value = 15510449
$ws.Range("N201").Value = 8368321
$ws.Range("O201").Value = 2605553
$ws.Range("F202").Value = 10071177
$ws.Range("M202").Value = 15673759
$ws.Range("N202").Value = 8381878
$ws.Range("O202").Value = 2635468
$ws.Range("F203").Value = 9851394
$ws.Range("B232").Value = 1361128
$ws.Range("F232").Value = 11588963
$ws.Range("G626").Value = 7009776
$ws.Range("H626").Value = 2635174
$ws.Range("I626").Value = 390.25
$ws.Range("M626").Value = 41276761
$ws.Range("M627").Value = 38624133
$ws.Range("M628").Value = 35818386
$ws.Range("M629").Value = 33636671
$ws.Range("M630").Value = 31534221
$ws.Range("M631").Value = 29240156
$ws.Range("M632").Value = 27537822
$ws.Range("M633").Value = 26108626
$ws.Range("M634").Value = 25082436
$ws.Range("M635").Value = 24559454
$ws.Range("F636").Value = 5707460
$ws.Range("M636").Value = 24210749
$ws.Range("F637").Value = 6576878
$ws.Range("M637").Value = 24256635
$ws.Range("R637").Value = 359372
$ws.Range("E638").Value = 140555
$ws.Range("J638").Value = 31.92
$ws.Range("O638").Value = 1481685
$ws.Range("P638").Value = 5156
$ws.Range("Q638").Value = 14408
$ws.Range("R638").Value = 360187
$ws.Range("C639").Value = 1870.1
$ws.Range("E639").Value = 118685
$ws.Range("J639").Value = 32.2
$ws.Range("O639").Value = 1461193
$ws.Range("P639").Value = 5171
$ws.Range("Q639").Value = 14444
$ws.Range("R639").Value = 359010
$ws.Range("A640").Value = 45016
$ws.Range("B640").Value = 986619
$ws.Range("C640").Value = 1852.7
$ws.Range("D640").Value = 410200
$ws.Range("E640").Value = 130478
$ws.Range("F640").Value = 7973206
$ws.Range("G640").Value = 6870522
$ws.Range("H640").Value = 2978415
$ws.Range("I640").Value = 441.25
$ws.Range("J640").Value = 32.73
$ws.Range("K640").Value = 13.72
$ws.Range("L640").Value = 413.68
$ws.Range("M640").Value = 25582401
$ws.Range("N640").Value = 4598297
$ws.Range("O640").Value = 1423214
$ws.Range("P640").Value = 5184
$ws.Range("Q640").Value = 14470
